# Added SVR parameter loading from pred_par structure and Excel files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the SVR parameters
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New parameter values
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Match the new active selection recorded after the edit
$ws.Range("I8").Select()
